# Lesson 2-6 - DP2.1.12 - Truth Tables
# Update the "DeMorgan's Law:" bullets on the Laws-of-Logic slide:
#   - the first law becomes "~(p v q) = ~p v ~q"
#   - the second (now redundant) law bullet is removed entirely

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)
$shp = $s.Shapes.Item(6)          # "Subtitle 2" placeholder
$tr = $shp.TextFrame.TextRange

# Paragraph 10 -> "(p v q) ^ r = (p ^ r) v (q ^ r)"  becomes the DeMorgan law
$tr.Paragraphs(10, 1).Text = "~(p v q) = ~p v ~q"

# Paragraph 11 -> "(p ^ q) ^ r = (p ^ r) ^ (q ^ r)" is deleted outright
$tr.Paragraphs(11, 1).Delete()
